$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Select row 2 (like a user clicking the row header) and delete the entire row.
# This is the row containing the 3.5mm audio jack part that isn't in the schematic.
$ws.Rows.Item(2).Select()
$ws.Rows.Item(2).Delete()

$excel.CalculateFullRebuild()
